# Remove the three OPQA-745 / OPQA-746 / OPQA-747 "evict user" test rows
# (rows 5-7) from the 1PAUTH sheet; the rows below shift up to fill the gap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5:7").Delete()

# The hyperlink table is not renumbered automatically when rows move, so
# rebuild it to match the new row layout (and drop the links that belonged
# to the deleted rows).
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A5"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-542", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-542")
$ws.Hyperlinks.Add($ws.Range("A2"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-539", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-539")
$ws.Hyperlinks.Add($ws.Range("A3"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-540", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-540")
$ws.Hyperlinks.Add($ws.Range("A4"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-541", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-541")
$ws.Hyperlinks.Add($ws.Range("A6"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-851", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-851")
$ws.Hyperlinks.Add($ws.Range("A7"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-852", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-852")
$ws.Hyperlinks.Add($ws.Range("A8"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-853", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-853")
$ws.Hyperlinks.Add($ws.Range("A9"), "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-854", "", "", "http://jira.bjz.apac.ime.reuters.com/browse/OPQA-854")

# Adding a hyperlink auto-applies the built-in "Hyperlink" cell style; the
# source workbook keeps these cells on the default style, so put it back.
$ws.Range("A2").Style = "Normal"
$ws.Range("A3").Style = "Normal"
$ws.Range("A4").Style = "Normal"
$ws.Range("A5").Style = "Normal"
$ws.Range("A6").Style = "Normal"
$ws.Range("A7").Style = "Normal"
$ws.Range("A8").Style = "Normal"
$ws.Range("A9").Style = "Normal"

# Selection previously covered L2:L12; with only 9 rows left it should
# cover L2:L9.
$ws.Range("L2:L9").Select()
